$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: "Ativação:" date value updated (01/01/2012 -> 01/01/2022) ---
# Use a text formula then convert to a plain value in place so Excel does not
# reinterpret the dd/mm/yyyy-looking text as a real date serial number and
# so the cell keeps its original style index (s="2" / s="3").
$ws.Cells.Item(8,2).Formula = "=""01/01/2022"""
$ws.Cells.Item(8,2).Copy()
$ws.Cells.Item(8,2).PasteSpecial(-4163)

$ws.Cells.Item(8,3).Formula = "=""01/01/2022"""
$ws.Cells.Item(8,3).Copy()
$ws.Cells.Item(8,3).PasteSpecial(-4163)

# --- Row 10: "Objetivos:" (Portuguese) text updated ---
$ws.Cells.Item(10,2).Value = "Proporcionar aos alunos uma visão atual dos processos industriais que utilizam a conversão química como rota de transformação da matéria prima em produto. Serão estudados os processos das indústrias de química de base e de transformação."
$ws.Cells.Item(10,3).Value = "Proporcionar aos alunos uma visão atual dos processos industriais que utilizam a conversão química como rota de transformação da matéria prima em produto. Serão estudados os processos das indústrias de química de base e de transformação."

# --- Row 11: new English objectives text added in B11/C11 (format copied from row 10) ---
$ws.Cells.Item(10,2).Copy()
$ws.Cells.Item(11,2).PasteSpecial(-4122)
$ws.Cells.Item(10,3).Copy()
$ws.Cells.Item(11,3).PasteSpecial(-4122)
$ws.Cells.Item(11,2).Value = "Provide students with a current view of industrial processes that use chemical conversion as a route to transform raw material into product. The processes of the basic chemical and transformation industries will be studied."
$ws.Cells.Item(11,3).Value = "Provide students with a current view of industrial processes that use chemical conversion as a route to transform raw material into product. The processes of the basic chemical and transformation industries will be studied."

# --- Row 14: "Programa resumido:" (Portuguese) text updated ---
$ws.Cells.Item(14,2).Value = "Introdução aos Processos Químicos Industriais; NPK / Fertilizantes; Ácido Sulfúrico; Cloro Álcalis; Papel e Celulose; Açúcar e álcool;  Processos Biotecnológicos;"
$ws.Cells.Item(14,3).Value = "Introdução aos Processos Químicos Industriais; NPK / Fertilizantes; Ácido Sulfúrico; Cloro Álcalis; Papel e Celulose; Açúcar e álcool;  Processos Biotecnológicos;"

# --- Row 15: new English short syllabus text added in B15/C15 (format copied from row 14) ---
$ws.Cells.Item(14,2).Copy()
$ws.Cells.Item(15,2).PasteSpecial(-4122)
$ws.Cells.Item(14,3).Copy()
$ws.Cells.Item(15,3).PasteSpecial(-4122)
$ws.Cells.Item(15,2).Value = "Introduction to Industrial Chemical Processes; NPK / Fertilizers; Sulfuric Acid; Chlorine Alkali; Paper and Cellulose; Sugar and alcohol; Biotechnological Processes."
$ws.Cells.Item(15,3).Value = "Introduction to Industrial Chemical Processes; NPK / Fertilizers; Sulfuric Acid; Chlorine Alkali; Paper and Cellulose; Sugar and alcohol; Biotechnological Processes."

# --- Row 16: "Programa:" (Portuguese) text updated ---
$ws.Cells.Item(16,2).Value = ".Introdução aos Processos Químicos Industriais; 2.NPK / Fertilizantes3.Ácido Sulfúrico; 4.Cloro Álcalis; 5.Papel e Celulose; 6.Açúcar e álcool; 7.Processos Biotecnológicos."
$ws.Cells.Item(16,3).Value = ".Introdução aos Processos Químicos Industriais; 2.NPK / Fertilizantes3.Ácido Sulfúrico; 4.Cloro Álcalis; 5.Papel e Celulose; 6.Açúcar e álcool; 7.Processos Biotecnológicos."

# --- Row 17: new English syllabus text added in B17/C17 (format copied from row 16) ---
$ws.Cells.Item(16,2).Copy()
$ws.Cells.Item(17,2).PasteSpecial(-4122)
$ws.Cells.Item(16,3).Copy()
$ws.Cells.Item(17,3).PasteSpecial(-4122)
$ws.Cells.Item(17,2).Value = "1. Introduction to Industrial Chemical Processes;2. NPK / Fertilizers3. Sulfuric Acid;4. Chlorine Alkali;5. Paper and Cellulose;6. Sugar and alcohol;7. Biotechnological Processes;"
$ws.Cells.Item(17,3).Value = "1. Introduction to Industrial Chemical Processes;2. NPK / Fertilizers3. Sulfuric Acid;4. Chlorine Alkali;5. Paper and Cellulose;6. Sugar and alcohol;7. Biotechnological Processes;"

# --- Row 19: "Método:" text updated ---
$ws.Cells.Item(19,2).Value = "Aulas expositivas, desenvolvimento de trabalhos e exercícios em sala e fora de sala de aula, discussão de casos práticos."
$ws.Cells.Item(19,3).Value = "Aulas expositivas, desenvolvimento de trabalhos e exercícios em sala e fora de sala de aula, discussão de casos práticos."

# --- Row 20: "Critério:" text updated ---
$ws.Cells.Item(20,2).Value = "Provas em sala, entrega de trabalhos e exercícios ou casos práticos elaborados fora de sala de aula."
$ws.Cells.Item(20,3).Value = "Provas em sala, entrega de trabalhos e exercícios ou casos práticos elaborados fora de sala de aula."

# --- Row 21: "Norma de recuperação:" text updated ---
$ws.Cells.Item(21,2).Value = "Frequência mínima de 70% e nota igual ou superior a 3,00 e inferior a 5,00 possibilita prova de recuperação."
$ws.Cells.Item(21,3).Value = "Frequência mínima de 70% e nota igual ou superior a 3,00 e inferior a 5,00 possibilita prova de recuperação."

# --- Row 22: "Bibliografia:" text updated ---
$biblio = "Ullmann" + [char]8217 + "s encyclopedia of industrial chemistry; Editorial advisory board, Giuseppe Bellussi et al.; 7th, completely revised edition; Weinheim ; New York : WileyVCH, 2011.Encyclopedia of Chemical Processing; Edited by Sunggyu Lee; New York : Taylor & Francis, 2006.Kirk, Raymond Eller. Encyclopedia of chemical technology / Herman F.Mark et al. New York: John Wiley, 1984.Manual econômico da indústria química - MEIQ / Centro de Pesquisas e Desenvolvimento; 8ed; Camaçari: CEPED, 2007.Shreve, R. Norris; BRINK JR., J. A. Indústrias de processos químicos. Tradução de Horácio Macedo; 4.ed. Rio de Janeiro: Editora Guanabara Koogan, 2008, c1997.T.W. Graham Solomons, Craig B. Fryhle Hoboken, NJ. Organic chemistry; John Wiley, 9th ed; c2008.Revistas:Brazilian Journal of Chemical Engineering, São Paulo, SP: Brazilian Society of Chemical Engineering, v. 11, n. 1, 1995-;"
$ws.Cells.Item(22,2).Value = $biblio
$ws.Cells.Item(22,3).Value = $biblio
